$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "[49.91766234016524, 50.14828405476167]"
$ws.Range("M3").Value = "[49.83296501771183, 50.13024019125025]"
$ws.Range("M4").Value = "[49.81732533343109, 50.14787995418201]"
$ws.Range("M5").Value = "[49.891966736767706, 50.14802422762657]"
$ws.Range("M6").Value = "[49.94111946031368, 50.25858049902765]"
$ws.Range("M7").Value = "[49.89163989941742, 50.17976938088707]"
$ws.Range("M8").Value = "[49.9059763492665, 50.18604143387028]"
$ws.Range("M9").Value = "[49.85151254633582, 50.13155946812045]"
$ws.Range("M10").Value = "[49.85641350336993, 50.1285091544404]"
$ws.Range("M11").Value = "[49.85407629522103, 50.167710200277305]"

$ws.Range("U2").Value = "[49.89337796154111, 50.06254316588849]"
$ws.Range("U3").Value = "[49.93621092807351, 50.095237999145695]"
$ws.Range("U4").Value = "[49.9732358155021, 50.14274559930923]"
$ws.Range("U5").Value = "[50.00427167366207, 50.15530320866451]"
$ws.Range("U6").Value = "[49.88333136767328, 50.06933995651849]"
$ws.Range("U7").Value = "[50.023557462964355, 50.17427073864197]"
$ws.Range("U8").Value = "[49.82701791426473, 49.983707749898784]"
$ws.Range("U9").Value = "[49.9241802157367, 50.09659367237862]"
$ws.Range("U10").Value = "[49.97597258939603, 50.14542998444228]"
$ws.Range("U11").Value = "[49.86363765574357, 50.02875981919525]"
